$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 122; this shifts the existing rows 122..247
# down to 123..248 and grows the sheet dimension to A1:T248.
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row 122 with the new price-observation record.
# Columns A,B,C,E,F,G,H,I,J,K,T are constant for every record in this sheet.
$ws.Cells.Item(122, 1).Value = 7
$ws.Cells.Item(122, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(122, 3).Value = "Ñuble"
$ws.Cells.Item(122, 4).Value = 44586
$ws.Cells.Item(122, 5).Value = 16
$ws.Cells.Item(122, 6).Value = "Fruta"
$ws.Cells.Item(122, 7).Value = 100101
$ws.Cells.Item(122, 8).Value = "Berries"
$ws.Cells.Item(122, 9).Value = 100112025
$ws.Cells.Item(122, 10).Value = "Frutilla"
$ws.Cells.Item(122, 11).Value = "Sin especificar"
$ws.Cells.Item(122, 12).Value = "Especial"
$ws.Cells.Item(122, 13).Value = 200
$ws.Cells.Item(122, 14).Value = 7000
$ws.Cells.Item(122, 15).Value = 7000
$ws.Cells.Item(122, 16).Value = 7000
$ws.Cells.Item(122, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(122, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(122, 19).Value = 1000
$ws.Cells.Item(122, 20).Value = 7
